# Recomputed line-loading results for the 380 kV case
# (Code/Results/Cases/Case_0_231/res_line/loading_percent.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 13  # columns C..O

$row2 = @(9.01820408189476,7.607254736013832,13.20875538047561,41.98635104916716,3.685679427710111,0,32.25314471626245,10.4744092580595,0,10.38439625577701,0,18.99769123883691,32.46513917364775)
$arr2 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr2[0,$i] = $row2[$i] }
$ws.Range("C2:O2").Value = $arr2

$row3 = @(9.026284143667457,7.619166327290476,13.22722596969196,41.76131505605715,3.689264724329762,0,32.12472334157368,10.49420576174005,0,10.40432801276011,0,18.40031508502705,32.31302885946887)
$arr3 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr3[0,$i] = $row3[$i] }
$ws.Range("C3:O3").Value = $arr3

$row4 = @(9.032606829518814,7.627175429834075,13.24046196189959,41.63325448983645,3.691583311878324,0,32.05343635976043,10.50761109073018,0,10.41772010662014,0,18.02485520896364,32.22733609144269)
$arr4 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr4[0,$i] = $row4[$i] }
$ws.Range("C4:O4").Value = $arr4

$row5 = @(9.035525839059773,7.630614320094867,13.24633228767249,41.58365123480633,3.69255773400047,0,32.02630787008246,10.51338847736726,0,10.42346796923212,0,17.86990355188765,32.19437749993968)
$arr5 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr5[0,$i] = $row5[$i] }
$ws.Range("C5:O5").Value = $arr5

$row6 = @(9.03603122401627,7.63119592959725,13.2473358334885,41.57557170401196,3.692721325709278,0,32.02191977967173,10.51436681616352,0,10.42443995059588,0,17.84406337566597,32.18902391359681)
$arr6 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr6[0,$i] = $row6[$i] }
$ws.Range("C6:O6").Value = $arr6

$row7 = @(9.032644809587213,7.627221098511666,13.24053920161131,41.63257501756477,3.691596333370152,0,32.05306269155859,10.50768773245974,0,10.41779644777081,0,18.02277304767602,32.22688362618261)
$arr7 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr7[0,$i] = $row7[$i] }
$ws.Range("C7:O7").Value = $arr7

$row8 = @(9.020707590515187,7.611217696648281,13.21473077934277,41.90668478266929,3.686891377876952,0,32.20730747324313,10.4809757574066,0,10.39102942067699,0,18.79364780656866,32.41110837230735)
$arr8 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr8[0,$i] = $row8[$i] }
$ws.Range("C8:O8").Value = $arr8

$row9 = @(9.008095190907783,7.585340521831191,13.17915328946066,42.52252545928049,3.678590008459284,0,32.56883389625045,10.43850311359026,0,10.34768134080416,0,20.2273683202997,32.83228033051572)
$arr9 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr9[0,$i] = $row9[$i] }
$ws.Range("C9:O9").Value = $arr9

$row10 = @(9.005399333771559,7.569669068581196,13.16217310825693,43.01989616275389,3.673048073627135,0,32.86899788218302,10.41332570842081,0,10.32138811905197,0,21.22223697909767,33.17643199513781)
$arr10 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr10[0,$i] = $row10[$i] }
$ws.Range("C10:O10").Value = $arr10

$row11 = @(9.005595722851368,7.563261710337153,13.15643464630134,43.25524313501465,3.670646410026934,0,33.01269722576998,10.40317753376363,0,10.31062878041591,0,21.66018057919901,33.34009884886685)
$arr11 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr11[0,$i] = $row11[$i] }
$ws.Range("C11:O11").Value = $arr11

$row12 = @(9.005874185416241,7.56093890835887,13.15454686261362,43.34560911248911,3.669754018713046,0,33.06810788966697,10.39952209487564,0,10.30672695497604,0,21.82377585682186,33.40305771350333)
$arr12 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr12[0,$i] = $row12[$i] }
$ws.Range("C12:O12").Value = $arr12

$row13 = @(9.005805145158529,7.561434564993388,13.15494075009878,43.32609281545938,3.669945453721603,0,33.05613052324649,10.40030102565641,0,10.3075596153684,0,21.78864458690801,33.38945539862658)
$arr13 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr13[0,$i] = $row13[$i] }
$ws.Range("C13:O13").Value = $arr13

$row14 = @(9.005614544916646,7.563068538649345,13.15627362253659,43.26265300682772,3.67057265099696,0,33.01723615636239,10.40287304306725,0,10.31030431932244,0,21.67368539489659,33.3452590760745)
$arr14 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr14[0,$i] = $row14[$i] }
$ws.Range("C14:O14").Value = $arr14

$row15 = @(9.005524359729998,7.564082869342791,13.15712718214455,43.22395457002656,3.670959047002838,0,32.99354075721487,10.40447288324399,0,10.31200798617252,0,21.60297336126124,33.31831416865841)
$arr15 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr15[0,$i] = $row15[$i] }
$ws.Range("C15:O15").Value = $arr15

$row16 = @(9.005415088211734,7.570102305024408,13.16258806578172,43.00469336455872,3.673207421731635,0,32.85974788798388,10.41401515957139,0,10.32211542192461,0,21.19330956972086,33.16587562754422)
$arr16 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr16[0,$i] = $row16[$i] }
$ws.Range("C16:O16").Value = $arr16

$row17 = @(9.00571212802164,7.573979693102517,13.16644657250394,42.87246869764549,3.674617233054374,0,32.77947897145218,10.42020316182844,0,10.32862356017369,0,20.93814219015166,33.07415235336734)
$arr17 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr17[0,$i] = $row17[$i] }
$ws.Range("C17:O17").Value = $arr17

$row18 = @(9.006016928227037,7.576277805706456,13.16885283058786,42.79727609406719,3.675439362867501,0,32.73398585690954,10.42388519847845,0,10.33247998486328,0,20.79000725568362,33.02206739528953)
$arr18 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr18[0,$i] = $row18[$i] }
$ws.Range("C18:O18").Value = $arr18

$row19 = @(9.006143148208508,7.577067583710449,13.16969966521625,42.77196664371775,3.675719656103634,0,32.71869967669892,10.42515298103309,0,10.33380514025136,0,20.73962067985785,33.00454892169395)
$arr19 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr19[0,$i] = $row19[$i] }
$ws.Range("C19:O19").Value = $arr19

$row20 = @(9.005666647555225,7.573559908625151,13.1660164816982,42.88645572307732,3.674465993300512,0,32.78795406793903,10.41953172435045,0,10.32791905275509,0,20.96544799484619,33.08384720834048)
$arr20 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr20[0,$i] = $row20[$i] }
$ws.Range("C20:O20").Value = $arr20

$row21 = @(9.005664994133594,7.562585793322602,13.15587438700849,43.28125352867688,3.670387965579219,0,33.02863365333585,10.40211249359345,0,10.30949345385297,0,21.70751365554066,33.35821429584388)
$arr21 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr21[0,$i] = $row21[$i] }
$ws.Range("C21:O21").Value = $arr21

$row22 = @(9.006853241747686,7.556016883257875,13.15090837898087,43.54650422214587,3.66782216897277,0,33.19171344840963,10.39182053544988,0,10.29845660139321,0,22.17935961385674,33.54323145197253)
$arr22 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr22[0,$i] = $row22[$i] }
$ws.Range("C22:O22").Value = $arr22

$row23 = @(9.006110413439092,7.559467714310561,13.15340684386442,43.4042947145695,3.669182518467133,0,33.10415728504513,10.39721365926365,0,10.30425528349495,0,21.92877110911181,33.44397646497094)
$arr23 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr23[0,$i] = $row23[$i] }
$ws.Range("C23:O23").Value = $arr23

$row24 = @(9.005686791768648,7.573749478404847,13.16621034029256,42.88012961397602,3.674534332619707,0,32.78412043535855,10.41983489350378,0,10.32823720292228,0,20.95310750188673,33.07946214291094)
$arr24 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr24[0,$i] = $row24[$i] }
$ws.Range("C24:O24").Value = $arr24

$row25 = @(9.010352202670299,7.591753207561138,13.18716891447912,42.34781854861411,3.680737427174831,0,32.46484992796635,10.4489336017068,0,10.34768134080416,0,19.84905939529497,32.71210494073937)
$arr25 = New-Object "object[,]" 1,$numCols
for ($i = 0; $i -lt $numCols; $i++) { $arr25[0,$i] = $row25[$i] }
$ws.Range("C25:O25").Value = $arr25

